$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (B9) was empty -> "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was a duplicate "Contact" / "No display for ContactDetail" row;
# turn it into the new "Jurisdiction" / "United States of America" row.
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# The old row 11 (second duplicate "Contact" row) is removed entirely,
# shifting the remaining rows (old 12-15, now 11-14) up by one.
$ws.Rows.Item(11).Delete()
